$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Video games rows -> LeagueOfLegendTV
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 3).Value = "LeagueOfLegendTV"
}

# Sport rows -> TennisTV
for ($r = 32; $r -le 61; $r++) {
    $ws.Cells.Item($r, 3).Value = "TennisTV"
}

# Music rows -> BestQualityMusic
for ($r = 62; $r -le 91; $r++) {
    $ws.Cells.Item($r, 3).Value = "BestQualityMusic"
}

# Widen the new column and update the active selection
$ws.Columns.Item(3).ColumnWidth = 20
[void]$ws.Range("C62").Select()
